$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.576.99'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '2.521.52'
$ws.Range('E3').Value = '  -2.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.97'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.89%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.566'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.522'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.67'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0805'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.108'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').Value = '2.905.98'
$ws.Range('E14').Value = '  -2.41%  '
$ws.Range('D15').Value = '2.560.74'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.805'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.44%  '
$ws.Range('D18').Value = '42.551.86'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.37%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.20'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '242.65'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.02'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.80%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.28'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.86%  '
$ws.Range('E28').Value = '  -2.93%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.09'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.62'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.95'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.74'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('E33').Value = '  +12.37%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.76%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0781'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.17'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.50%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.99'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.07%  '
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.15'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '21.67'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -10.17%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0300'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.25'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').Value = '2.000.06'
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').Value = '2.759.06'
$ws.Range('E48').Value = '  -2.53%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.189'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.34%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '72.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.54%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.87%  '
